# Portugal Segunda Liga workbook update
# - Swap mismatched fixture data (columns B:AC) between several pairs of adjacent
#   rows whose odds/result data had been attributed to the wrong fixture.
# - Remove three fixtures (rows 370-372) that no longer belong in the sheet,
#   shifting the remaining rows up and renumbering the sequential id column (A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of rows whose B:AC content (everything except the sequential id in column A)
# needs to be swapped between the two rows.
$pairs = @(
    @(5, 6),
    @(41, 42),
    @(44, 45),
    @(71, 72),
    @(74, 75),
    @(81, 82),
    @(90, 91),
    @(254, 255),
    @(312, 313),
    @(320, 321)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B$r1`:AC$r1")
    $rng2 = $ws.Range("B$r2`:AC$r2")

    $v1 = $rng1.Value()
    $v2 = $rng2.Value()

    $rng1.Value = $v2
    $rng2.Value = $v1
}

# Remove the three fixtures that were dropped from the source data.
$ws.Rows("370:372").Delete()

# After the row shift, re-sequence the id column (A) so that it keeps matching
# "row number - 2", as it does throughout the rest of the sheet.
$lastRow = $ws.UsedRange.Rows.Count + $ws.UsedRange.Row - 1
for ($r = 370; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
